$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column width updates (D, E, F)
# Note: the engine's ColumnWidth setter stores (set value + 5/6) as the
# OOXML <col width>, matching observed round-trip behaviour, so we
# back the desired stored widths (11, 22, 18) out by that offset.
$ws.Columns.Item(4).ColumnWidth = 11 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 18 - (5/6)

# Row 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 344.284604629486
$ws.Range("F2").Value = 0

# Row 3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3120.1145
$ws.Range("F3").Value = 0

# Row 4
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 250.631825420901
$ws.Range("F4").Value = 0

# Row 13
$ws.Range("C13").Value = 130
$ws.Range("E13").Value = 130

# Row 14
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 240
$ws.Range("F14").Value = 0

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1638
$ws.Range("F15").Value = 0

# Row 16
$ws.Range("C16").Value = 17085.89
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 17085.89
$ws.Range("F16").Value = 0

# Row 17
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 342
$ws.Range("F17").Value = 0

# Row 18
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 1200
$ws.Range("F18").Value = 0

# Row 19
$ws.Range("C19").Value = 27181.31093005039
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 27181.31093005039
$ws.Range("F19").Value = 0
